$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update Price (D) and Volume(1h) (E) columns for rows with new market data ---
$ws.Range("D2").Value = "70.590.29"
$ws.Range("E2").Value = "  -0.94%  "
$ws.Range("D3").Value = "3.795.34"
$ws.Range("E3").Value = "  -1.90%  "
$ws.Range("E4").Value = "  +0.06%  "
$ws.Range("D5").Value = "'705.23"
$ws.Range("E5").Value = "  +2.10%  "
$ws.Range("D6").Value = "'168.84"
$ws.Range("E6").Value = "  -2.55%  "
$ws.Range("D7").Value = "3.790.55"
$ws.Range("E7").Value = "  -2.12%  "
$ws.Range("E8").Value = "  +0.13%  "
$ws.Range("E9").Value = "  -1.29%  "
$ws.Range("D10").Value = "'0.159"
$ws.Range("E10").Value = "  -2.36%  "
$ws.Range("D11").Value = "'7.32"
$ws.Range("E11").Value = "  -1.10%  "
$ws.Range("E12").Value = "  -1.77%  "
$ws.Range("D14").Value = "'36.08"
$ws.Range("E14").Value = "  -1.60%  "
$ws.Range("D15").Value = "4.438.95"
$ws.Range("E15").Value = "  -1.66%  "
$ws.Range("D16").Value = "3.850.19"
$ws.Range("E16").Value = "  -0.07%  "
$ws.Range("D17").Value = "70.629.54"
$ws.Range("E17").Value = "  -0.91%  "
$ws.Range("E18").Value = "  -0.10%  "
$ws.Range("D19").Value = "'7.13"
$ws.Range("E19").Value = "  -1.80%  "
$ws.Range("E20").Value = "  -3.41%  "
$ws.Range("D21").Value = "'490.06"
$ws.Range("E21").Value = "  -0.28%  "
$ws.Range("D22").Value = "'10.49"
$ws.Range("E22").Value = "  -5.35%  "
$ws.Range("D23").Value = "'0.724"
$ws.Range("E23").Value = "  +0.09%  "
$ws.Range("D24").Value = "'84.99"
$ws.Range("E24").Value = "  +0.15%  "
$ws.Range("E25").Value = "  -3.50%  "
$ws.Range("D26").Value = "'12.02"
$ws.Range("E26").Value = "  -2.98%  "
$ws.Range("D27").Value = "'10.41"
$ws.Range("E27").Value = "  -1.42%  "
$ws.Range("D28").Value = "3.951.24"
$ws.Range("E28").Value = "  -1.61%  "
$ws.Range("E29").Value = "  -0.11%  "
$ws.Range("D30").Value = "'2.05"
$ws.Range("E30").Value = "  -4.73%  "
$ws.Range("D31").Value = "'3.07"
$ws.Range("E31").Value = "  -1.54%  "
$ws.Range("D32").Value = "'7.31"
$ws.Range("E32").Value = "  -4.45%  "
$ws.Range("E33").Value = "  -4.51%  "
$ws.Range("D34").Value = "'29.02"
$ws.Range("E34").Value = "  -2.85%  "
$ws.Range("D35").Value = "'0.172"
$ws.Range("E35").Value = "  -4.46%  "
$ws.Range("D38").Value = "'9.01"
$ws.Range("E38").Value = "  -3.27%  "
$ws.Range("E39").Value = "  -3.35%  "
$ws.Range("E40").Value = "  +1.51%  "
$ws.Range("D41").Value = "'2.29"
$ws.Range("E41").Value = "  -4.40%  "
$ws.Range("D42").Value = "'5.89"
$ws.Range("E42").Value = "  -3.14%  "
$ws.Range("E43").Value = "  -5.75%  "
$ws.Range("E44").Value = "  -0.06%  "
$ws.Range("E45").Value = "  +0.17%  "
$ws.Range("D46").Value = "'164.25"
$ws.Range("E46").Value = "  +0.18%  "
$ws.Range("D47").Value = "'0.000306"
$ws.Range("E47").Value = "  -1.20%  "
$ws.Range("D48").Value = "'48.70"
$ws.Range("E48").Value = "  -0.04%  "
$ws.Range("D49").Value = "'420.49"
$ws.Range("E49").Value = "  +1.18%  "
$ws.Range("D50").Value = "'8.66"
$ws.Range("E50").Value = "  -0.53%  "
$ws.Range("D51").Value = "'0.292"
$ws.Range("E51").Value = "  -3.74%  "

# --- Rows 36 and 37 swapped rank order: RenzoRestakedETH <-> Binance-PegBSC-USD ---
$ws.Range("B36").Value = "Binance-PegBSC-USD"
$ws.Range("C36").Value = "https://coinranking.com/coin/i5jggxiwp+binance-pegbsc-usd-bsc-usd"
$ws.Range("D36").Value = "'1.00"
$ws.Range("E36").Value = "  +0.18%  "

$ws.Range("B37").Value = "RenzoRestakedETH"
$ws.Range("C37").Value = "https://coinranking.com/coin/lKlJ_MC5M+renzorestakedeth-ezeth"
$ws.Range("D37").Value = "3.767.17"
$ws.Range("E37").Value = "  -1.28%  "

